$wb = $excel.ActiveWorkbook

# --- survey sheet: add "hideInContents" column (S) -------------------------
$survey = $wb.Worksheets.Item("survey")
$survey.Range("S1").Value = "hideInContents"
$survey.Range("S3").Value = $true
$survey.Range("S6").Value = $true
$survey.Range("S8").Value = $true
$survey.Range("S11").Value = $true
$survey.Range("S14").Value = $true
$survey.Range("S17").Value = $true

# --- section1 sheet: add "hideInContents" column (T) ------------------------
$section1 = $wb.Worksheets.Item("section1")
$section1.Range("T1").Value = "hideInContents"
$section1.Range("T2").Value = $true

# --- section2 sheet: add "hideInContents" column (T) ------------------------
$section2 = $wb.Worksheets.Item("section2")
$section2.Range("T1").Value = "hideInContents"
$section2.Range("T2").Value = $true

# --- view/selection state ----------------------------------------------------
$section1.Activate() | Out-Null
$section1.Range("T2").Select() | Out-Null

$section2.Activate() | Out-Null
$section2.Range("T3").Select() | Out-Null

$survey.Activate() | Out-Null
$survey.Range("S18").Select() | Out-Null
